# Generate Report for handback
#
# Refreshes the "Latest Handoff Datetime" and "Latest Handback DateTime"
# columns (D and G) for the first data row (row 2 - the
# "0c38d88b-cef0-44dd-8447-a71ab02f33f1.md" source file) on both the
# "zh-cn" and "de-de" localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-08 11:11:57"
$zhcn.Range("G2").Value = "2016-01-08 11:12:41"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-08 11:12:07"
$dede.Range("G2").Value = "2016-01-08 11:12:57"
